# ----------------------------------------------------------------------------
# Refresh the "cryptos" worksheet with the latest scraped Price / Volume(1h)
# figures (and a handful of rank re-orderings among same-value rows) that came
# in with the "Updated cryptos list ... with GitHub Actions" commit.
#
# All of the data cells in this sheet are stored as *text* (coinranking.com
# formats its prices with dotted thousands separators like "44.552.31", and
# the % change column keeps literal padding spaces), so every write below is a
# plain string assignment. For cells whose new text happens to *look* like a
# plain number (e.g. "6.79"), Excel's COM layer would otherwise auto-coerce
# the cell to a numeric type on assignment; we avoid that by writing it the way
# a user typing into the grid would force text (a leading apostrophe) and then
# resetting the cell style back to Normal so no stray number-format style is
# left behind.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    if ($Text -match "^-?\d+(\.\d+)?$") {
        # Numeric-looking text: force text entry so it is not coerced to a number,
        # then strip the number-format style that the forced-text entry adds.
        $Range.Formula = "'" + $Text
        $Range.Style = "Normal"
    } else {
        $Range.Value = $Text
    }
}

Set-TextValue $ws.Range("D2") "44.552.31"
Set-TextValue $ws.Range("E2") "  +2.55%  "
Set-TextValue $ws.Range("D3") "2.368.84"
Set-TextValue $ws.Range("E3") "  -0.11%  "
Set-TextValue $ws.Range("D5") "0.673"
Set-TextValue $ws.Range("E5") "  +3.41%  "
Set-TextValue $ws.Range("D6") "239.60"
Set-TextValue $ws.Range("E6") "  +2.84%  "
Set-TextValue $ws.Range("D7") "73.34"
Set-TextValue $ws.Range("E7") "  +5.56%  "
Set-TextValue $ws.Range("D9") "0.553"
Set-TextValue $ws.Range("E9") "  +20.31%  "
Set-TextValue $ws.Range("D10") "0.102"
Set-TextValue $ws.Range("E10") "  +7.01%  "
Set-TextValue $ws.Range("D11") "29.93"
Set-TextValue $ws.Range("E11") "  +12.07%  "
Set-TextValue $ws.Range("B12") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D12") "2.721.61"
Set-TextValue $ws.Range("E12") "  -0.28%  "
Set-TextValue $ws.Range("B13") "TRON"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D13") "0.107"
Set-TextValue $ws.Range("E13") "  +1.61%  "
Set-TextValue $ws.Range("D14") "16.94"
Set-TextValue $ws.Range("E14") "  +7.32%  "
Set-TextValue $ws.Range("D15") "6.79"
Set-TextValue $ws.Range("D16") "0.902"
Set-TextValue $ws.Range("E16") "  +6.90%  "
Set-TextValue $ws.Range("D17") "2.367.19"
Set-TextValue $ws.Range("E17") "  -0.27%  "
Set-TextValue $ws.Range("D18") "44.527.65"
Set-TextValue $ws.Range("E19") "  +4.95%  "
Set-TextValue $ws.Range("D20") "77.61"
Set-TextValue $ws.Range("D21") "6.50"
Set-TextValue $ws.Range("E21") "  +3.82%  "
Set-TextValue $ws.Range("D22") "255.36"
Set-TextValue $ws.Range("E22") "  +2.58%  "
Set-TextValue $ws.Range("D24") "3.79"
Set-TextValue $ws.Range("E24") "  -2.80%  "
Set-TextValue $ws.Range("E25") "  +3.20%  "
Set-TextValue $ws.Range("E26") "  +4.11%  "
Set-TextValue $ws.Range("D27") "2.24"
Set-TextValue $ws.Range("E27") "  -1.27%  "
Set-TextValue $ws.Range("D28") "22.53"
Set-TextValue $ws.Range("E28") "  +0.22%  "
Set-TextValue $ws.Range("D29") "1.61"
Set-TextValue $ws.Range("E29") "  +5.69%  "
Set-TextValue $ws.Range("D30") "174.33"
Set-TextValue $ws.Range("E30") "  -0.28%  "
Set-TextValue $ws.Range("E31") "  +2.07%  "
Set-TextValue $ws.Range("E32") "  +5.53%  "
Set-TextValue $ws.Range("D33") "0.0748"
Set-TextValue $ws.Range("E33") "  +8.10%  "
Set-TextValue $ws.Range("E34") "  +4.42%  "
Set-TextValue $ws.Range("D35") "5.24"
Set-TextValue $ws.Range("E35") "  +4.08%  "
Set-TextValue $ws.Range("D36") "3.93"
Set-TextValue $ws.Range("E36") "  +7.90%  "
Set-TextValue $ws.Range("E37") "  -3.17%  "
Set-TextValue $ws.Range("D38") "6.53"
Set-TextValue $ws.Range("E39") "  +6.94%  "
Set-TextValue $ws.Range("D40") "20.18"
Set-TextValue $ws.Range("E40") "  +9.95%  "
Set-TextValue $ws.Range("E41") "  +0.16%  "
Set-TextValue $ws.Range("D42") "8.87"
Set-TextValue $ws.Range("E42") "  -1.33%  "
Set-TextValue $ws.Range("E43") "  +3.35%  "
Set-TextValue $ws.Range("B44") "Cronos"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D44") "0.0986"
Set-TextValue $ws.Range("E44") "  +3.59%  "
Set-TextValue $ws.Range("B45") "ARBITRUM"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D45") "1.18"
Set-TextValue $ws.Range("E45") "  +0.61%  "
Set-TextValue $ws.Range("B46") "Algorand"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D46") "0.184"
Set-TextValue $ws.Range("E46") "  +12.17%  "
Set-TextValue $ws.Range("B47") "FTXToken"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D47") "4.48"
Set-TextValue $ws.Range("E47") "  +1.93%  "
Set-TextValue $ws.Range("B48") "Aave"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D48") "98.91"
Set-TextValue $ws.Range("E48") "  -0.20%  "
Set-TextValue $ws.Range("D49") "2.36"
Set-TextValue $ws.Range("E49") "  +3.37%  "
Set-TextValue $ws.Range("D50") "1.445.23"
Set-TextValue $ws.Range("E50") "  -0.11%  "
Set-TextValue $ws.Range("D51") "2.595.83"
Set-TextValue $ws.Range("E51") "  -0.18%  "
